$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab (hyphen instead of space before "10")
$ws.Name = "Результаты ТОП-10"

# Update the leaderboard names (column B) and scores (column C) for rows 2-9
$ws.Range("B2").Value = "Смельчак"
$ws.Range("C2").Value = 160.0

$ws.Range("B3").Value = "Нюша хрюша"
$ws.Range("C3").Value = 156.0

$ws.Range("B4").Value = "Отвинта"
$ws.Range("C4").Value = 143.0

$ws.Range("B5").Value = "Копатыч"
$ws.Range("C5").Value = 136.0

$ws.Range("B6").Value = ""
$ws.Range("C6").Value = 112.0

$ws.Range("B7").Value = "Винни"
$ws.Range("C7").Value = 90.0

$ws.Range("B8").Value = "Совунья"
$ws.Range("C8").Value = 66.0

$ws.Range("B9").Value = "Крош"
$ws.Range("C9").Value = 34.0

# Add a new 10th row
$ws.Range("A10").Value = 9.0
$ws.Range("B10").Value = ""
$ws.Range("C10").Value = 0.0
